# Updated symbol list on Wed Dec 21 08:31:13 UTC 2022 with GitHub Actions
# Applies the price/volume-label refresh described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $rng = $ws.Range($CellRef)
    # Prefix with an apostrophe so Excel stores the value as text rather than
    # re-parsing numeric-looking strings into floating point numbers (which
    # would lose formatting such as trailing zeros / exact decimal digits).
    $rng.Value = "'" + $NewValue
    # Reset the cell style back to Normal so we don't leave a stray
    # quote-prefix style applied to the cell (keeps styling identical to
    # the original workbook).
    $rng.Style = "Normal"
}

# Column D (Price) updates
Set-TextCell "D3"  "22.44"
Set-TextCell "D4"  "5.388"
Set-TextCell "D5"  "0.05674"
Set-TextCell "D6"  "3.407"
Set-TextCell "D7"  "6.311"
Set-TextCell "D8"  "0.8081"
Set-TextCell "D9"  "0.9244"
Set-TextCell "D10" "0.1404"
Set-TextCell "D11" "0.07422"
Set-TextCell "D12" "0.03079"
Set-TextCell "D13" "0.03022"
Set-TextCell "D14" "0.09368"
Set-TextCell "D17" "0.04754"
Set-TextCell "D19" "0.0005855"
Set-TextCell "D20" "0.006447"
Set-TextCell "D21" "0.004976"
Set-TextCell "D22" "0.001008"
Set-TextCell "D25" "2.153"
Set-TextCell "D26" "0.3256"
Set-TextCell "D27" "0.1299"
Set-TextCell "D41" "0.003009"
Set-TextCell "D43" "0.002711"
Set-TextCell "D44" "0.007537"
Set-TextCell "D45" "0.00005806"
Set-TextCell "D47" "0.4003"
Set-TextCell "D48" "0.2118"

# Column E (Volume(1h)) updates
Set-TextCell "E19" "18OneONE"
Set-TextCell "E41" "40KickTokenKICKWorstin24h"
